$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (in order) that receive new values for rows 2..25
$cols = @(2,3,5,6,7,8,9,10,11,13)  # B, C, E, F, G, H, I, J, K, M

# New values per row (row 2 first), one array per row matching $cols order
$data = @(
    @(0.2719510070050148,0.04192620222919174,0.395888443509115,2.482631370127592,0.9367928068124769,0.9788314240462199,0.8155962367436231,0.07284871740412902,0.2930086416530742,0.3782657008355841),
    @(0.2423707676603328,0.0366973643492372,0.3882734587432708,2.468090838625741,0.9393248354716519,0.9847937332924346,0.8215949119766819,0.07357365049769471,0.2598701233779792,0.3597497190560048),
    @(0.2242431601743533,0.03348659013961708,0.3837728690030815,2.460501578136672,0.9414815795330043,0.9888943725769366,0.8257062547252794,0.07404177443787008,0.2395528692035072,0.3485475407146978),
    @(0.2168650648345363,0.03217810945902499,0.3819829294790225,2.457745385411258,0.9425115979926915,0.9906759742264057,0.8274892137811207,0.0742383261036228,0.2312812165130111,0.3440246333372556),
    @(0.2156404935295768,0.0319608333203405,0.3816883762887642,2.457308041306263,0.9426917516976516,0.9909784845956864,0.8277917657077758,0.07427131292451894,0.2299081935343281,0.3432761535042133),
    @(0.224143619541195,0.03346894377315834,0.3837485506489031,2.460463044890076,0.941494859169282,0.9889179522682241,0.825729864994937,0.07404440176876559,0.2394412828905104,0.3484863726575398),
    @(0.2617446958370806,0.04012334712983545,0.3932265063789089,2.477339860636206,0.9375407284091892,0.9807959423528416,0.8175756584754907,0.07309389871424732,0.2815764177312587,0.3718468733052092),
    @(0.335746628666783,0.05317142656188878,0.4131998845196065,2.521068117812021,0.934576832791663,0.9683593627620866,0.8049872332842867,0.0714124925046491,0.3644333248334988,0.4189756465396641),
    @(0.3902709570770924,0.06275917020026611,0.4287197050049301,2.559699994414956,0.9353388256778459,0.9613526366291296,0.7978195821231608,0.07028830954720888,0.4254448658978163,0.4544040956145494),
    @(0.4151080169842487,0.06712167302778482,0.4359637336215414,2.578692453103386,0.93632803354933,0.9586282779710444,0.7950122893103853,0.06980099685388597,0.4532299083733733,0.4706958140047846),
    @(0.424517786905426,0.06877380141416722,0.4387332801131123,2.586088693113467,0.9367953434229861,0.9576632630787145,0.7940145527006877,0.06961992691064278,0.4637556109091463,0.4768901642909142),
    @(0.4224910289935906,0.06841797961574514,0.438135636069191,2.584486694891808,0.9366905711615772,0.9578681315351503,0.7942265255990719,0.06965876943870342,0.4614885343081028,0.4755549893515649),
    @(0.4158820758086961,0.06725759174636892,0.4361910573606522,2.579296852368032,0.9363646193855999,0.958547549711966,0.7949288949261017,0.06978603067173728,0.4540957835275208,0.4712049261707563),
    @(0.4118344793920699,0.06654683935987293,0.4350033815555889,2.576144524199961,0.936177048578898,0.9589723932264604,0.7953676276983117,0.06986443311898283,0.449568036942253,0.4685436419612827),
    @(0.3886484462787507,0.0624740903686245,0.4282499872166809,2.558487356975334,0.935287135944364,0.9615400021670553,0.7980121789112715,0.07032064133755966,0.4236296387938125,0.4533429044925725),
    @(0.3744329557913488,0.05997584485901086,0.4241540736461289,2.548018750591595,0.9349060075085021,0.963233781818289,0.7997507273531994,0.070606679826279,0.4077248995537559,0.4440624911996522),
    @(0.3662597701652999,0.05853901323857258,0.4218155339109586,2.542130999940795,0.9347472556221135,0.9642515788010257,0.8007933535554912,0.0707734683768062,0.3985798102070248,0.4387411370150716),
    @(0.3634930275929094,0.05805254271453464,0.4210267215540711,2.540160433627619,0.934703879899601,0.9646036705229477,0.8011536913846129,0.07083032938904799,0.3954839463094402,0.4369422555350155),
    @(0.3759458914019262,0.06024177771685402,0.4245882984224565,2.549119331792269,0.9349403192756398,0.9630489653335701,0.7995612398341692,0.07057599592697761,0.4094176885710397,0.4450487010882256),
    @(0.4178231656383673,0.06759842189555343,0.4367615115725911,2.580815691076467,0.9364578404488242,0.9583461790631986,0.7947208180374616,0.0697485569161218,0.4562671043754563,0.4724819670893012),
    @(0.4452185484530844,0.07240724977391721,0.4448712012992999,2.602721474684856,0.9379902085860437,0.9556610955330314,0.7919381470356726,0.06922797399075353,0.4869097300708063,0.4905570162260275),
    @(0.4305948305447487,0.06984061097659833,0.4405288571563659,2.590920955440836,0.9371227870410479,0.957058612256759,0.7933884194287657,0.06950397038273604,0.4705530986680913,0.4808967301374238),
    @(0.375261894779868,0.06012155120598095,0.4243919347623617,2.548621351635546,0.9349246189148488,0.9631323837433001,0.7996467728575354,0.07058986082464092,0.4086523825081656,0.4446027917328479),
    @(0.3156993931665966,0.04964149280932872,0.4076480838395824,2.508097809685182,0.9348638249147427,0.9713498394891218,0.8080277960892843,0.07184783019465257,0.3419940992525312,0.4060849867553031)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $c = $cols[$j]
        $ws.Cells.Item($r, $c).Value = $rowValues[$j]
    }
}
